# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (and before the
#    "总计" summary sheet), mirroring the per-quarter holding sheets already
#    in the workbook (968013 / Schroder Asian Income fund).
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q1"

# Reuse the header/formatting already used by the sibling quarter sheets so
# the new sheet's style indices line up with the existing ones instead of
# minting new cellXfs entries.
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Range("A2:H2").Copy($newSheet.Range("A2:H2"))
$templateSheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))

$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# The source data cells (fund code / size / weight figures) are stored as
# text in the source system, not numbers - force text formatting before
# assigning so "968013" etc. round-trip as strings instead of being
# auto-coerced to numeric cells.
$newSheet.Range("B2:G2").NumberFormat = "@"

$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "968013"
$newSheet.Cells.Item(2,3).Value = "施罗德亚洲高息股债基金M"
$newSheet.Cells.Item(2,4).Value = "297.64"
$newSheet.Cells.Item(2,5).Value = "57.54"
$newSheet.Cells.Item(2,6).Value = "2.08"
$newSheet.Cells.Item(2,7).Value = "6.1909"
$newSheet.Cells.Item(2,8).Value = 1

# Drop the transient "@" text-format style picked up above: stamp the plain
# (unformatted) look of an untouched cell back over the data row so the
# saved styles line up with the rest of the workbook, which never carries an
# explicit number format on this row.
$blankCell = $newSheet.Cells.Item(100,100)
$blankCell.Copy()
$newSheet.Range("B2:G2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Update the "总计" roll-up sheet: push the existing quarters down one row
#    and insert the new 2022-Q1 summary row at the top of the data block.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Grab the formatting of the current last index cell (A4) and stamp it onto
# the row that is about to appear below it (A5) so the new row carries the
# same style id as the rest of the index column.
$summary.Cells.Item(4,1).Copy($summary.Cells.Item(5,1))

for ($r = 4; $r -ge 2; $r--) {
    $dest = $r + 1
    $summary.Cells.Item($dest,1).Value = $dest - 2
    $summary.Cells.Item($dest,2).Value = $summary.Cells.Item($r,2).Value2
    $summary.Cells.Item($dest,3).Value = $summary.Cells.Item($r,3).Value2
    $summary.Cells.Item($dest,4).Value = $summary.Cells.Item($r,4).Value2
}

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q1"
$summary.Cells.Item(2,3).Value = 1
$summary.Cells.Item(2,4).Value = 6.19
